$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "light" column (col C) for the new
# "text" field, pushing light/background/dark/complementary/vector right.
$ws.Columns.Item(3).Insert()

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "color"
$ws.Range("B1").Value = "base"
$ws.Range("C1").Value = "text"
$ws.Range("D1").Value = "light"
$ws.Range("E1").Value = "background"
$ws.Range("F1").Value = "dark"
$ws.Range("G1").Value = "complementary"
$ws.Range("H1").Value = "vector"

# --- Row 2: "autumn" (was "green") --------------------------------------
$ws.Range("A2").Value = "autumn"
$ws.Range("B2").Value = "#00313C"
$ws.Range("C2").Value = "#00313C"
$ws.Range("D2").Value = "#1F4D5B"
$ws.Range("E2").Value = "#F9FAF7"
$ws.Range("F2").Value = "#1D2834"
$ws.Range("G2").Value = "#81402C"
$ws.Range("H2").Value = "#00313C, #183C41, #304746, #48514A, #856949, #AB7743, #D2863C, #FB9637, #E57630, #D0562B, #BC3626, #A61922, #8D0422, #70002A, #520036"
$ws.Range("H2").Font.Name = "Helvetica"

# --- Row 3: "rainbow" (was "gray") --------------------------------------
$ws.Range("A3").Value = "rainbow"
$ws.Range("B3").Value = "#C65D44"
$ws.Range("C3").Value = "#213C47"
$ws.Range("D3").Value = "#F4A261"
$ws.Range("E3").Value = "#FCFAF7"
$ws.Range("F3").Value = "#652D1F"
$ws.Range("G3").Value = "#44ACC5"
$ws.Range("H3").Value = "#1B3037, #213C47, #26525B, #28847D, #45A289, #7CAE7F, #B2B974, #E9C46A, #ECBA67, #EFB065, #F2A662, #F0935C, #E87653, #D8674B, #C65D44"
$ws.Range("H3").Font.Name = "Helvetica"

# --- Selection matches the committed state: whole row 3 selected -------
$ws.Range("A3:XFD3").Select()
